# Update sha256 for passwords
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$hash = "055ef921716dabc901c32a25c6b158db41087db41d7c997aa713b4777dedcc65"

# Replace every password in column B (rows 2-6) with the new sha256 hash
$ws.Range("B2").Value = $hash
$ws.Range("B3").Value = $hash
$ws.Range("B4").Value = $hash
$ws.Range("B5").Value = $hash
$ws.Range("B6").Value = $hash

# Row 6 user record changed: username + id
$ws.Range("A6").Value = "reutisa"
$ws.Range("C6").Value = 207555555

# Column B now needs to be wide enough to show the full hash
$ws.Columns.Item(2).ColumnWidth = 86.9

# Update the active selection to match the saved state
$ws.Range("B12").Select()
